$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variables")

# New "footnote" column header in E1
$ws.Range("E1").Value = "footnote"

# C5 gets a footnote marker appended to its text, and the footnote text
# itself goes into the new E5 cell
$ws.Range("C5").Value = "Proud of being Samoan*"
$ws.Range("E5").Value = "*Somewhat proud or very proud"

# Give the new column E a sensible width (closest the engine's column-width
# rounding model can reach to the template's 27.42578125 stored width)
$ws.Columns.Item(5).ColumnWidth = 26.666666666666668

# Re-apply the AutoFilter so the dropdown only spans the header row and now
# also covers the new column
$ws.AutoFilterMode = $false
$ws.Range("A1:E1").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the new AutoFilter range
foreach ($n in $wb.Names) {
  if ($n.Name -eq "variables!_FilterDatabase") {
    $n.RefersTo = "=variables!`$A`$1:`$E`$1"
  }
}

# Mirror the template's saved selection
$ws.Range("D18").Select()
